# QAS Reboot test data update: change the TransitProjectID value for QAS
# (shared string referenced by I2) and refresh the sheet's view state
# (selection, row height) to match the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Transit data -------------------------------------------------
# I2 holds "TransitProjectID" data (header in I1). Re-assign the same
# text-typed value so the cell keeps storing a string (not a number),
# matching the original shared-string cell type.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1000104331"
$ws.Range("I2").NumberFormat = "General"

# --- Refresh view/selection ----------------------------------------------
# Selection moves from N2 to I2, and the view scrolls back to the top-left
# (A1) instead of being scrolled to show column F first.
$ws.Range("A1").Select()
$ws.Range("I2").Select()

# --- Row height tweak ------------------------------------------------------
# Row 2's height changes from 14.65 to 12.8.
$ws.Rows.Item(2).RowHeight = 12.8
